# Generate Report for Handoff
#
# Moves the localization status report from "In Translation" to
# "Ready for handoff" and refreshes the associated timestamps, widening
# the Status columns so the new (longer) text fits.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# zh-cn / de-de status columns (E, F)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-16 14:39:39"

# Widen the now-longer Status columns to fit "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- "zh-cn" sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 14:39:34"
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- "de-de" sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 14:39:39"
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
